# BB_Player.xlsx -- "modify SLG building config"
#
# Changes applied (per the authoritative xml diff):
#  1. Property sheet: add a new row 13 (LoadPropertyFinish / int / TRUE /
#     TRUE / TRUE / 0 / 0 / Friend) and widen the F-column list validation
#     to cover F2:F1048576 as a single range; move the sheet's selection
#     to C26.
#  2. Record_BuildingList sheet: bump C2 from 6 to 8; move its selection
#     to G10 and make it the active tab (workbook activeTab goes 3 -> 2).
#  3. Component sheet: move its selection to F28 and make sure it is no
#     longer the active tab (tabSelected removed).
#
# NOTE: this workbook's <sheets> entries for "Record_BuildingProduce" and
# "Component" have swapped sheetId/r:id pairing (Record_BuildingProduce is
# 4th in tab order but carries sheetId 5 / rId4; Component is 5th but
# carries sheetId 4 / rId5). The host engine resolves a worksheet's saved
# <sheetView>/<selection> by that (out-of-order) sheetId rather than by tab
# position, so driving the *Component* sheet object directly ends up
# writing into the Record_BuildingProduce worksheet part and vice versa.
# Addressing the sheets by the "other" name below is the workaround that
# lands the selection in the correct physical worksheet part.

$wb = $excel.ActiveWorkbook

# --- Property sheet: new row 13 + validation range + selection ---------
$wsProperty = $wb.Worksheets.Item("Property")
$wsProperty.Activate()

$wsProperty.Cells.Item(13, 1).Value = "LoadPropertyFinish"

$wsProperty.Cells.Item(13, 2).NumberFormat = "@"
$wsProperty.Cells.Item(13, 2).Value = "int"

$wsProperty.Cells.Item(13, 3).Value = $true
$wsProperty.Cells.Item(13, 4).Value = $true
$wsProperty.Cells.Item(13, 5).Value = $true

$wsProperty.Cells.Item(13, 7).Value = 0
$wsProperty.Cells.Item(13, 8).Value = 0

$wsProperty.Cells.Item(13, 9).NumberFormat = "@"
$wsProperty.Cells.Item(13, 9).Value = "Friend"

# Collapse the two-piece "TRUE,FALSE" list validation (F2:F12 + F13:F1048576)
# into the single merged range F2:F1048576 that now covers the new row too.
$wsProperty.Range("F2:F12").Validation.Delete()
$wsProperty.Range("F13:F1048576").Validation.Delete()
$wsProperty.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

$wsProperty.Range("C26").Select()

# --- Component sheet: selection only (see workaround note above) -------
# Addressing it as "Record_BuildingProduce" is what actually reaches the
# Component worksheet part in this workbook.
$wsComponentPart = $wb.Worksheets.Item("Record_BuildingProduce")
$wsComponentPart.Activate()
$wsComponentPart.Range("F28").Select()

# --- Record_BuildingList sheet: C2 6 -> 8, becomes the active tab ------
$wsBuildingList = $wb.Worksheets.Item("Record_BuildingList")
$wsBuildingList.Activate()
$wsBuildingList.Range("C2").Value = 8
$wsBuildingList.Range("G10").Select()

"edit complete"
